# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) figures.
#
# Each new value is written with a leading apostrophe escape (`') so Excel
# stores it as literal text -- exactly like the original inlineStr cells --
# instead of auto-converting number/date look-alike strings (e.g. "222.96"
# or multi-dot figures like "27.239.05") into numeric or date values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'27.239.05"
$ws.Range("E2").Value = "`'  -2.65%  "
$ws.Range("D3").Value = "`'1.704.08"
$ws.Range("E3").Value = "`'  -1.72%  "
$ws.Range("D4").Value = "`'1.003"
$ws.Range("E4").Value = "`'  +0.18%  "
$ws.Range("D5").Value = "`'222.96"
$ws.Range("D6").Value = "`'0.5297"
$ws.Range("E6").Value = "`'  -2.52%  "
$ws.Range("D7").Value = "`'1.004"
$ws.Range("E7").Value = "`'  +0.17%  "
$ws.Range("D8").Value = "`'0.2649"
$ws.Range("E8").Value = "`'  -4.43%  "
$ws.Range("D9").Value = "`'0.06576"
$ws.Range("E9").Value = "`'  -2.38%  "
$ws.Range("D10").Value = "`'20.76"
$ws.Range("E10").Value = "`'  -4.43%  "
$ws.Range("D11").Value = "`'0.07647"
$ws.Range("E11").Value = "`'  -2.23%  "
$ws.Range("D12").Value = "`'4.573"
$ws.Range("E12").Value = "`'  -3.00%  "
$ws.Range("D13").Value = "`'1.720.88"
$ws.Range("E13").Value = "`'  -0.70%  "
$ws.Range("D14").Value = "`'1.938.64"
$ws.Range("E14").Value = "`'  -1.66%  "
$ws.Range("D15").Value = "`'0.5720"
$ws.Range("E15").Value = "`'  -4.63%  "
$ws.Range("D16").Value = "`'0.0₅8171"
$ws.Range("D17").Value = "`'67.43"
$ws.Range("E17").Value = "`'  -2.32%  "
$ws.Range("D18").Value = "`'27.215.23"
$ws.Range("E18").Value = "`'  -2.48%  "
$ws.Range("D19").Value = "`'216.13"
$ws.Range("E19").Value = "`'  -1.34%  "
$ws.Range("E20").Value = "`'  +0.18%  "
$ws.Range("D21").Value = "`'4.658"
$ws.Range("E21").Value = "`'  -3.45%  "
$ws.Range("D22").Value = "`'10.43"
$ws.Range("E22").Value = "`'  -4.53%  "
$ws.Range("D23").Value = "`'5.961"
$ws.Range("E23").Value = "`'  -4.53%  "
$ws.Range("D24").Value = "`'1.004"
$ws.Range("E24").Value = "`'  +0.20%  "
$ws.Range("D25").Value = "`'142.16"
$ws.Range("E25").Value = "`'  -2.88%  "
$ws.Range("D26").Value = "`'1.742"
$ws.Range("E26").Value = "`'  +6.30%  "
$ws.Range("D27").Value = "`'0.1218"
$ws.Range("E27").Value = "`'  -2.24%  "
$ws.Range("D28").Value = "`'7.247"
$ws.Range("E28").Value = "`'  -2.76%  "
$ws.Range("D29").Value = "`'16.27"
$ws.Range("E29").Value = "`'  -3.89%  "
$ws.Range("D30").Value = "`'0.05361"
$ws.Range("E30").Value = "`'  -4.46%  "
$ws.Range("D31").Value = "`'1.289"
$ws.Range("E31").Value = "`'  -2.31%  "
$ws.Range("D32").Value = "`'3.504"
$ws.Range("E32").Value = "`'  -5.89%  "
$ws.Range("D33").Value = "`'3.408"
$ws.Range("E33").Value = "`'  -3.65%  "
$ws.Range("D34").Value = "`'1.632"
$ws.Range("E34").Value = "`'  -0.54%  "
$ws.Range("D35").Value = "`'2.873"
$ws.Range("E35").Value = "`'  +0.67%  "
$ws.Range("D36").Value = "`'2.421"
$ws.Range("E36").Value = "`'  -0.96%  "
$ws.Range("D37").Value = "`'0.9445"
$ws.Range("E37").Value = "`'  -4.08%  "
$ws.Range("D38").Value = "`'0.5847"
$ws.Range("E38").Value = "`'  -0.99%  "
$ws.Range("D39").Value = "`'0.01629"
$ws.Range("E39").Value = "`'  -2.46%  "
$ws.Range("D40").Value = "`'5.859"
$ws.Range("E40").Value = "`'  -1.37%  "
$ws.Range("D41").Value = "`'1.004"
$ws.Range("E41").Value = "`'  +0.22%  "
$ws.Range("D42").Value = "`'1.038.30"
$ws.Range("E42").Value = "`'  -0.88%  "
$ws.Range("D43").Value = "`'0.8379"
$ws.Range("E43").Value = "`'  -0.51%  "
$ws.Range("D44").Value = "`'100.93"
$ws.Range("E44").Value = "`'  -1.44%  "
$ws.Range("D45").Value = "`'1.846.35"
$ws.Range("E45").Value = "`'  -1.59%  "
$ws.Range("D46").Value = "`'0.0₈115"
$ws.Range("E46").Value = "`'  -2.97%  "
$ws.Range("D47").Value = "`'57.93"
$ws.Range("D48").Value = "`'0.4489"
$ws.Range("E48").Value = "`'  +1.68%  "
$ws.Range("D49").Value = "`'1.001"
$ws.Range("E49").Value = "`'  -0.25%  "
$ws.Range("D50").Value = "`'0.06603"
$ws.Range("E50").Value = "`'  +11.43%  "
$ws.Range("D51").Value = "`'8.051"
$ws.Range("E51").Value = "`'  -2.39%  "
